# Edit corresponding to the commit "Add files via upload":
#   - Slide 4: merge the "анализ " / "существующих моделей " runs into one run.
#   - Slide 5: split "Можно рассмотреть итеративный модель, ..." into
#     "Можно " / "рассмотреть " / "итеративную " / "модель, ..." runs,
#     fixing "итеративный" -> "итеративную" along the way.

$p = $ppt.ActivePresentation

# --- Slide 4 ("Прямоугольник 2"): combine "анализ " and "существующих моделей " ---
$s4 = $p.Slides.Item(4)
$shape4 = $s4.Shapes.Item(1)
$tr4 = $shape4.TextFrame.TextRange
$bullet1 = $tr4.Paragraphs(2)
$combined = $tr4.Characters($bullet1.Start, 28)
$combined.Text = "анализ существующих моделей "

# --- Slide 5: last bullet, split single run into four ---
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(1)
$tr5 = $shape5.TextFrame.TextRange
$bullet4 = $tr5.Paragraphs(6)
$start5 = $bullet4.Start

# "итеративный " (chars 19-30 of the bullet) -> "итеративную "
$wordRange = $tr5.Characters($start5 + 18, 12)
$wordRange.Text = "итеративную "

# Split "Можно рассмотреть " into "Можно " + "рассмотреть "
$secondRange = $tr5.Characters($start5 + 6, 12)
$secondRange.Text = "рассмотреть "
